$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.058134008642893
$ws.Range("D2").Value = 1.062455705162352
$ws.Range("E2").Value = 1.053963137347904
$ws.Range("F2").Value = 1.071201952948032
$ws.Range("I2").Value = 1.047987611823277
$ws.Range("J2").Value = 1.063126810246375
$ws.Range("K2").Value = 1.06517700479328
$ws.Range("L2").Value = 1.056707657155979
$ws.Range("M2").Value = 1.073899755022633

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.059862937996661
$ws.Range("D3").Value = 1.063832423060787
$ws.Range("E3").Value = 1.055481907089025
$ws.Range("F3").Value = 1.072715604502015
$ws.Range("I3").Value = 1.048502448399858
$ws.Range("J3").Value = 1.06450474203489
$ws.Range("K3").Value = 1.066367129738078
$ws.Range("L3").Value = 1.058037771461454
$ws.Range("M3").Value = 1.075228193111517

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.060978667832774
$ws.Range("D4").Value = 1.06472058198062
$ws.Range("E4").Value = 1.056461518941099
$ws.Range("F4").Value = 1.073692580818812
$ws.Range("I4").Value = 1.048832764465002
$ws.Range("J4").Value = 1.065393040585171
$ws.Range("K4").Value = 1.067134021105517
$ws.Range("L4").Value = 1.0588948165296
$ws.Range("M4").Value = 1.076084818791365

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.061447016811284
$ws.Range("D5").Value = 1.065093335774026
$ws.Range("E5").Value = 1.056872610763556
$ws.Range("F5").Value = 1.074102725050798
$ws.Range("I5").Value = 1.04897095987755
$ws.Range("J5").Value = 1.065765699724994
$ws.Range("K5").Value = 1.067455666623572
$ws.Range("L5").Value = 1.059254262085241
$ws.Range("M5").Value = 1.076444245917837

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.061525613827975
$ws.Range("D6").Value = 1.065155886251875
$ws.Range("E6").Value = 1.056941592053395
$ws.Range("F6").Value = 1.074171556597444
$ws.Range("I6").Value = 1.048994124384609
$ws.Range("J6").Value = 1.065828225353812
$ws.Range("K6").Value = 1.067509628307764
$ws.Range("L6").Value = 1.059314564724134
$ws.Range("M6").Value = 1.07650455468072

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.060984928679717
$ws.Range("D7").Value = 1.064725565188528
$ws.Range("E7").Value = 1.056467014847937
$ws.Range("F7").Value = 1.073698063435705
$ws.Range("I7").Value = 1.048834613662381
$ws.Range("J7").Value = 1.065398023129205
$ws.Range("K7").Value = 1.067138321903632
$ws.Range("L7").Value = 1.058899622806071
$ws.Range("M7").Value = 1.07608962420306

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.058718939341728
$ws.Range("D8").Value = 1.062921532772319
$ws.Range("E8").Value = 1.054477069542585
$ws.Range("F8").Value = 1.071714014488179
$ws.Range("I8").Value = 1.048162189738353
$ws.Range("J8").Value = 1.06359318258187
$ws.Range("K8").Value = 1.065579881671759
$ws.Range("L8").Value = 1.05715793412514
$ws.Range("M8").Value = 1.074349326934912

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.054702296365188
$ws.Range("D9").Value = 1.059721659471585
$ws.Range("E9").Value = 1.05094593520868
$ws.Range("F9").Value = 1.068198517476201
$ws.Range("I9").Value = 1.046955483554464
$ws.Range("J9").Value = 1.0603868823516
$ws.Range("K9").Value = 1.062808747583995
$ws.Range("L9").Value = 1.05406051604946
$ws.Range("M9").Value = 1.071259517233207

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.052007657544746
$ws.Range("D10").Value = 1.057573641600839
$ws.Range("E10").Value = 1.048574471348
$ws.Range("F10").Value = 1.065841100597632
$ws.Range("I10").Value = 1.046136025715468
$ws.Range("J10").Value = 1.05823113736911
$ws.Range("K10").Value = 1.060943890426828
$ws.Range("L10").Value = 1.051975736321832
$ws.Range("M10").Value = 1.069183362327295

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.050836632389139
$ws.Range("D11").Value = 1.056639867987466
$ws.Range("E11").Value = 1.047543293001747
$ws.Range("F11").Value = 1.064816884663705
$ws.Range("I11").Value = 1.045777563910668
$ws.Range("J11").Value = 1.057293183857207
$ws.Range("K11").Value = 1.060132105555585
$ws.Range("L11").Value = 1.051068124609696
$ws.Range("M11").Value = 1.06828034676959

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.050401007657316
$ws.Range("D12").Value = 1.056292458381121
$ws.Range("E12").Value = 1.047159602349567
$ws.Range("F12").Value = 1.064435914150284
$ws.Range("I12").Value = 1.04564386334023
$ws.Range("J12").Value = 1.056944094849261
$ws.Range("K12").Value = 1.05982991538922
$ws.Range("L12").Value = 1.050730247642173
$ws.Range("M12").Value = 1.067944308074646

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.050494480528784
$ws.Range("D13").Value = 1.056367004684393
$ws.Range("E13").Value = 1.04724193565691
$ws.Range("F13").Value = 1.064517657886601
$ws.Range("I13").Value = 1.045672567643766
$ws.Range("J13").Value = 1.057019007152128
$ws.Range("K13").Value = 1.059894766147056
$ws.Range("L13").Value = 1.050802757563159
$ws.Range("M13").Value = 1.068016417723265

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.050800636959164
$ws.Range("D14").Value = 1.056611162595421
$ws.Range("E14").Value = 1.047511590675233
$ws.Range("F14").Value = 1.064785404441023
$ws.Range("I14").Value = 1.045766523475363
$ws.Range("J14").Value = 1.057264342240737
$ws.Range("K14").Value = 1.060107139901145
$ws.Range("L14").Value = 1.051040210973337
$ws.Range("M14").Value = 1.068252582402512

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.050989182841146
$ws.Range("D15").Value = 1.056761521133533
$ws.Range("E15").Value = 1.04767764541029
$ws.Range("F15").Value = 1.064950301116838
$ws.Range("I15").Value = 1.045824339441109
$ws.Range("J15").Value = 1.0574154092293
$ws.Range("K15").Value = 1.060237902923338
$ws.Range("L15").Value = 1.05118641407533
$ws.Range("M15").Value = 1.068398008889643

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.052085283956068
$ws.Range("D16").Value = 1.057635534659153
$ws.Range("E16").Value = 1.048642814883365
$ws.Range("F16").Value = 1.065909000752556
$ws.Range("I16").Value = 1.046159738599958
$ws.Range("J16").Value = 1.05829329013411
$ws.Range("K16").Value = 1.060997674462921
$ws.Range("L16").Value = 1.052035867241294
$ws.Range("M16").Value = 1.069243206504865

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.052771694054881
$ws.Range("D17").Value = 1.058182788839463
$ws.Range("E17").Value = 1.049247072507124
$ws.Range("F17").Value = 1.066509436867205
$ws.Range("I17").Value = 1.046369149006053
$ws.Range("J17").Value = 1.058842745875994
$ws.Range("K17").Value = 1.06147310177123
$ws.Range("L17").Value = 1.052567387383254
$ws.Range("M17").Value = 1.069772289052646

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.053171658699874
$ws.Range("D18").Value = 1.058501639983026
$ws.Range("E18").Value = 1.049599110156375
$ws.Range("F18").Value = 1.06685933027904
$ws.Range("I18").Value = 1.046490944632908
$ws.Range("J18").Value = 1.05916280078532
$ws.Range("K18").Value = 1.061749997474784
$ws.Range("L18").Value = 1.052876943468188
$ws.Range("M18").Value = 1.070080506347617

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.053307967872826
$ws.Range("D19").Value = 1.05861030052012
$ws.Range("E19").Value = 1.049719075893093
$ws.Range("F19").Value = 1.066978579213263
$ws.Range("I19").Value = 1.046532414632706
$ws.Range("J19").Value = 1.059271858188925
$ws.Range("K19").Value = 1.061844342083201
$ws.Range("L19").Value = 1.052982414723736
$ws.Range("M19").Value = 1.070185535142712

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.052698090953404
$ws.Range("D20").Value = 1.058124110250785
$ws.Range("E20").Value = 1.049182284487373
$ws.Range("F20").Value = 1.066445050076735
$ws.Range("I20").Value = 1.046346717491067
$ws.Range("J20").Value = 1.058783839404239
$ws.Range("K20").Value = 1.061422135730705
$ws.Range("L20").Value = 1.052510409090183
$ws.Range("M20").Value = 1.0697155636918

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.050710499720432
$ws.Range("D21").Value = 1.056539279868397
$ws.Range("E21").Value = 1.047432202484727
$ws.Range("F21").Value = 1.064706574532445
$ws.Range("I21").Value = 1.045738871111203
$ws.Range("J21").Value = 1.057192116376705
$ws.Range("K21").Value = 1.060044619366205
$ws.Range("L21").Value = 1.050970307676626
$ws.Range("M21").Value = 1.068183054974101

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.049457029924785
$ws.Range("D22").Value = 1.055539561229508
$ws.Range("E22").Value = 1.046328000218802
$ws.Range("F22").Value = 1.063610445817473
$ws.Range("I22").Value = 1.045353497663065
$ws.Range("J22").Value = 1.056187329673078
$ws.Range("K22").Value = 1.05917471186443
$ws.Range("L22").Value = 1.049997640403431
$ws.Range("M22").Value = 1.06721592058968

$ws.Range("B23").Value = 1.019999999999999
$ws.Range("C23").Value = 1.050121883685967
$ws.Range("D23").Value = 1.056069845768155
$ws.Range("E23").Value = 1.046913729730894
$ws.Range("F23").Value = 1.064191821222834
$ws.Range("I23").Value = 1.045558096534513
$ws.Range("J23").Value = 1.056720371004389
$ws.Range("K23").Value = 1.059636231546261
$ws.Range("L23").Value = 1.050513686628482
$ws.Range("M23").Value = 1.067728961232717

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.052731350283306
$ws.Range("D24").Value = 1.058150625669727
$ws.Range("E24").Value = 1.049211560688679
$ws.Range("F24").Value = 1.066474144717925
$ws.Range("I24").Value = 1.046356854408302
$ws.Range("J24").Value = 1.05881045804198
$ws.Range("K24").Value = 1.061445166366618
$ws.Range("L24").Value = 1.052536156582179
$ws.Range("M24").Value = 1.069741196638185

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.055743598831266
$ws.Range("D25").Value = 1.06055145536722
$ws.Range("E25").Value = 1.051861818507636
$ws.Range("F25").Value = 1.069109726669904
$ws.Range("I25").Value = 1.047270063626385
$ws.Range("J25").Value = 1.061218940711611
$ws.Range("K25").Value = 1.063528176960801
$ws.Range("L25").Value = 1.054864711110113
$ws.Range("M25").Value = 1.072061125483575

